# Applies:
#  1. Split ", February 2012" into ", " / "Social Informatics, " / "February 2012"
#  2. Prefix the exercise paragraph with an italic "Exercise 1: " run, and move the
#     "_GoBack" bookmark to sit right after that new run.
#  3. (handled implicitly by 2 above) Remove the "_GoBack" bookmark from its old spot
#     at the end of the long paragraph -- re-adding a bookmark with the same name
#     relocates the (singleton) "_GoBack" bookmark instead of creating a duplicate.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: "University of Tartu, February 2012" -> "University of Tartu, Social
#          Informatics, February 2012", split across three runs.
# ---------------------------------------------------------------------------
$para2 = $d.Paragraphs.Item(2)
$commaSpace = $para2.Range.Start + [int]"University of Tartu".Length

$insertionPoint = $d.Range($commaSpace + 2, $commaSpace + 2)
$insertionPoint.InsertBefore("Social Informatics, ")

# Force the newly-touched text to stay split into distinct runs (rather than being
# re-coalesced with its neighbours at save time) by toggling a character property
# on/off -- this leaves no visible formatting trace but "dirties" the run.
$rComma = $d.Range($commaSpace, $commaSpace + 2)
$rComma.Bold = 1
$rComma.Bold = 0

$rSocial = $d.Range($commaSpace + 2, $commaSpace + 2 + "Social Informatics, ".Length)
$rSocial.Bold = 1
$rSocial.Bold = 0

# ---------------------------------------------------------------------------
# Part 2: Prefix "Read "What is Social..." with an italic "Exercise 1: " run and
#          move the _GoBack bookmark to sit between the two runs.
# ---------------------------------------------------------------------------
$readFind = $d.Content
$readFind.Find.Execute("Read") | Out-Null

# Isolate "Read" into its own run first (preserves its italic/lang formatting)
$readFind.Bold = 1
$readFind.Bold = 0

$readStart = $readFind.Start

# Grow the text via FormattedText so the new prefix inherits "Read"'s formatting
$fmt = $readFind.FormattedText
$fmt.Text = "Exercise 1: Read"

$splitPos = $readStart + "Exercise 1: ".Length

# Re-add the singleton "_GoBack" bookmark here; this both creates it at the new
# location (splitting "Exercise 1: " and "Read ..." into separate runs) and
# removes it from its old location further down in the document.
$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos)) | Out-Null
